$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add new header "Save" in H1, matching the style of the other header cells (e.g. G1)
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("H1").Value = "Save"

# Fill in the "Save" column values for rows 2-6
$ws.Range("H2").Value = 0
$ws.Range("H3").Value = 1
$ws.Range("H4").Value = 0
$ws.Range("H5").Value = 0
$ws.Range("H6").Value = 1
